# Adds traceability for factory
# Rebuilds the tail of the traceability sheet so a new "Factory Traceability"
# table (with its own Java/C#/Test-case rows) is inserted between the main
# Bridge table and the existing Iterator Traceability table. Also gives the
# whole sheet visible thin-box gridlines (previously borderless) and resizes
# the columns to fit the widened content.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlContinuous = 1
$xlLineStyleNone = -4142
$xlEdgeLeft = 7
$xlEdgeTop = 8
$xlEdgeBottom = 9
$xlEdgeRight = 10

function Set-CellBorders($cell, $left, $right, $top, $bottom) {
    if ($left)   { $cell.Borders.Item($xlEdgeLeft).LineStyle   = $xlContinuous } else { $cell.Borders.Item($xlEdgeLeft).LineStyle   = $xlLineStyleNone }
    if ($right)  { $cell.Borders.Item($xlEdgeRight).LineStyle  = $xlContinuous } else { $cell.Borders.Item($xlEdgeRight).LineStyle  = $xlLineStyleNone }
    if ($top)    { $cell.Borders.Item($xlEdgeTop).LineStyle    = $xlContinuous } else { $cell.Borders.Item($xlEdgeTop).LineStyle    = $xlLineStyleNone }
    if ($bottom) { $cell.Borders.Item($xlEdgeBottom).LineStyle = $xlContinuous } else { $cell.Borders.Item($xlEdgeBottom).LineStyle = $xlLineStyleNone }
}

# Apply the same border shape to every cell (A,B,C) of a row - matches how
# the workbook's table rows are styled (each cell individually boxed, not
# just the outer edge of the 3-cell range).
function Set-RowBoxStyle($rowNum, $left, $right, $top, $bottom) {
    foreach ($col in @("A","B","C")) {
        Set-CellBorders ($ws.Range("$col$rowNum")) $left $right $top $bottom
    }
}

function Set-FullBoxRows($firstRow, $lastRow) {
    for ($r = $firstRow; $r -le $lastRow; $r++) {
        Set-RowBoxStyle $r $true $true $true $true
    }
}

function Set-Row3($rowNum, $a, $b, $c) {
    $ws.Range("A$rowNum").Value = $a
    $ws.Range("B$rowNum").Value = $b
    $ws.Range("C$rowNum").Value = $c
}

# ---------------------------------------------------------------------
# 1. Remove the old trailing content (J13/C13/T11 row, the Iterator
#    Traceability block, and the old Factory Traceability block) so we can
#    rebuild it in the new order.
# ---------------------------------------------------------------------
$ws.Rows("15:19").Delete()

# 2. Make room for the new layout: blank / Factory header / Factory table
#    header / 3 factory data rows / old J13 row / blank / Iterator header /
#    Iterator data row = 10 rows.
$ws.Rows("15:24").Insert()

# ---------------------------------------------------------------------
# 3. Fill in the new rows with their final content.
# ---------------------------------------------------------------------
# Row 15: blank spacer
Set-Row3 15 "" "" ""

# Row 16: "Factory Traceability" banner (merged)
Set-Row3 16 "Factory Traceability" "" ""
$ws.Range("A16:C16").Merge()

# Row 17: column headers for the factory table
Set-Row3 17 "Java Code" "C# Code" "Test Case(s)"

# Rows 18-20: factory traceability data
Set-Row3 18 "J14" "C14" "T12"
Set-Row3 19 "J15" "C15" "T12, T13"
Set-Row3 20 "J16" "C16" "T12, T13"

# Row 21: the old J13/C13/T11 row, now trailing the factory table
Set-Row3 21 "J13" "C13" "T11"

# Row 22: blank spacer
Set-Row3 22 "" "" ""

# Row 23: "Iterator Traceability" banner (merged)
Set-Row3 23 "Iterator Traceability" "" ""
$ws.Range("A23:C23").Merge()

# Row 24: Iterator traceability detail row (unchanged content, shifted down)
Set-Row3 24 "Iterator.java: lines 5-52" " Iterator.cs all (modifying to support generic type instead of just integers)" " iteratorTest"

# ---------------------------------------------------------------------
# 4. Borders: the whole sheet now gets visible thin-box gridlines (it had
#    none before), plus the open-box look around each banner row.
# ---------------------------------------------------------------------
Set-FullBoxRows 1 13
Set-RowBoxStyle 14 $true $true $true $false   # top edge only (closes box from above)
Set-RowBoxStyle 15 $false $false $true $true  # spacer: horizontal rules only
Set-RowBoxStyle 16 $true $true $false $true   # bottom edge only (closes box from below)
Set-FullBoxRows 17 20
Set-RowBoxStyle 21 $true $true $true $false
Set-RowBoxStyle 22 $false $false $true $true
Set-RowBoxStyle 23 $true $true $false $true
Set-FullBoxRows 24 24

# ---------------------------------------------------------------------
# 5. Column widths / sheet cosmetics matching the resized table.
# ---------------------------------------------------------------------
$ws.Columns("A").ColumnWidth = 20.375
$ws.Columns("B").ColumnWidth = 102.25
$ws.Columns("C").ColumnWidth = 24.25

$ws.Range("B27").Select()
